# ---------------------------------------------------------------------------
# teacher_import_template.xlsx update:
#   * "Teachers" sheet: drop the "email" column, shift remaining headers
#     left, and replace the 2 sample rows with 5 new sample teacher rows.
#   * Add a new "Instructions" sheet right after "Teachers".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Teachers")

# Drop the now-unused last column (previously "emergency_contact" in column J;
# the sheet now only needs columns A-I).
$ws.Columns.Item(10).Delete()

# --- Header row (row 1) ------------------------------------------------------
$headers = @("first_name", "last_name", "phone_number", "subjects", "classes_assigned", "qualification", "experience_years", "address", "emergency_contact")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Sample data rows (rows 2-6) ---------------------------------------------
# Columns: first_name, last_name, phone_number, subjects, classes_assigned,
#          qualification, experience_years, address, emergency_contact
$rows = @(
    @("Priya",  "Sharma", "9876601001", "Mathematics,Science",    "Class 7,Class 8",  "M.Sc Mathematics",       8,  "15 Koramangala, Bangalore",  "9876602001"),
    @("Ravi",   "Verma",  "9876601002", "English,Hindi",          "Class 9,Class 10", "M.A English Literature", 5,  "32 Banjara Hills, Hyderabad","9876602002"),
    @("Meena",  "Nair",   "9876601003", "Social Studies,History", "Class 7,Class 9",  "M.A History",            12, "67 Andheri West, Mumbai",    "9876602003"),
    @("Suresh", "Patel",  "9876601004", "Physics,Chemistry",      "Class 8,Class 10", "M.Sc Physics",           6,  "89 Satellite, Ahmedabad",    "9876602004"),
    @("Kavita", "Joshi",  "9876601005", "Biology,Mathematics",    "Class 7,Class 8",  "M.Sc Zoology",           10, "44 Sadashiv Peth, Pune",     "9876602005")
)

# Phone-number-like columns (C = phone_number, I = emergency_contact) must be
# stored as text, not auto-converted to numbers.
$ws.Range("C2:C6").NumberFormat = "@"
$ws.Range("I2:I6").NumberFormat = "@"

for ($r = 0; $r -lt $rows.Count; $r++) {
    $rowData = $rows[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $rowData[$c]
    }
}

# --- Add the "Instructions" sheet, placed right after "Teachers" -----------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$ws2.Name = "Instructions"

$instructions = @(
    "Instructions",
    "1. Fill in teacher details in the Teachers sheet",
    "2. Required fields: first_name, last_name, phone_number",
    "3. Phone numbers: 10 digits (e.g., 9876543210)",
    '4. Subjects: Comma-separated (e.g., "Mathematics,Science")',
    '5. Classes: Comma-separated (e.g., "Class 7,Class 8")',
    "6. Phone number must be unique",
    "7. Email is auto-generated from phone number",
    "8. Upload this file in Admin Panel -> Import Data",
    "",
    "Note: This template has 5 sample teachers ready to import!",
    "",
    "Login Credentials:",
    "- Teachers log in via mobile app using OTP",
    "- Use the phone number from the import",
    "- No password needed for mobile login",
    "- OTP will be sent to the phone number"
)

for ($i = 0; $i -lt $instructions.Count; $i++) {
    $ws2.Cells.Item($i + 1, 1).Value = $instructions[$i]
}

# Give the "Instructions" title cell the same bold/bordered header look used
# for the Teachers header row.
$ws.Range("A1").Copy() | Out-Null
$ws2.Range("A1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

Write-Host "Edit complete"
